$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "35.631.36"
    "E2" = "  -2.74%  "
    "D3" = "1.981.85"
    "E3" = "  -3.70%  "
    "E4" = "  -0.15%  "
    "D5" = "246.32"
    "E5" = "  +1.17%  "
    "E6" = "  -4.81%  "
    "D7" = "58.53"
    "E7" = "  +7.80%  "
    "E8" = "  -0.03%  "
    "E9" = "  -0.24%  "
    "D10" = "0.361"
    "E10" = "  -0.59%  "
    "D11" = "0.0736"
    "E11" = "  -1.83%  "
    "E12" = "  -2.73%  "
    "D13" = "0.941"
    "E13" = "  +1.17%  "
    "D14" = "14.57"
    "E14" = "  -0.97%  "
    "D15" = "2.270.85"
    "E15" = "  -3.80%  "
    "E16" = "  -2.63%  "
    "D17" = "1.997.95"
    "E17" = "  -3.36%  "
    "D18" = "18.20"
    "E18" = "  +7.31%  "
    "D19" = "35.530.76"
    "E19" = "  -2.83%  "
    "D20" = "71.43"
    "E20" = "  -0.76%  "
    "D21" = "0.0₃0848"
    "E21" = "  -1.61%  "
    "D22" = "5.23"
    "E22" = "  -0.43%  "
    "D23" = "232.87"
    "E23" = "  -2.13%  "
    "B24" = "Dai"
    "C24" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D24" = "1.00"
    "E24" = "  +0.08%  "
    "B25" = "PancakeSwap"
    "C25" = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    "D25" = "2.63"
    "E25" = "  +23.01%  "
    "D26" = "2.27"
    "D27" = "164.90"
    "E27" = "  +0.48%  "
    "D28" = "9.13"
    "E28" = "  -2.19%  "
    "D29" = "19.18"
    "E29" = "  -4.53%  "
    "E30" = "  -2.42%  "
    "D31" = "4.86"
    "E31" = "  -4.19%  "
    "E32" = "  -6.00%  "
    "D33" = "0.0955"
    "E33" = "  +15.95%  "
    "D34" = "0.0595"
    "E34" = "  -0.20%  "
    "E35" = "  +9.94%  "
    "D36" = "4.35"
    "E36" = "  -3.22%  "
    "E37" = "  -0.10%  "
    "E38" = "  -3.92%  "
    "D39" = "5.36"
    "E39" = "  +9.42%  "
    "E40" = "  -2.13%  "
    "E41" = "  -0.59%  "
    "D42" = "0.0213"
    "E42" = "  -1.23%  "
    "D43" = "7.88"
    "E43" = "  +3.65%  "
    "D44" = "93.48"
    "E44" = "  -0.77%  "
    "E45" = "  -1.34%  "
    "D46" = "16.18"
    "E46" = "  +1.29%  "
    "D47" = "0.0898"
    "E47" = "  -0.87%  "
    "D48" = "1.378.28"
    "E48" = "  -1.95%  "
    "D49" = "2.90"
    "E49" = "  -0.46%  "
    "D50" = "47.27"
    "E50" = "  +4.15%  "
    "E51" = "  +0.05%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates"
